$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Range("K2").Value = 912
$ws.Range("K3").Value = 852
$ws.Range("K4").Value = 197
$ws.Range("K5").Value = 51
$ws.Range("J6").Value = 11057
$ws.Range("K6").Value = 1183
$ws.Range("J7").Value = 29250
$ws.Range("K7").Value = 3195

$ws = $wb.Worksheets.Item('Austin')
$ws.Range("K2").Value = 50
$ws.Range("K6").Value = 65
$ws.Range("K7").Value = 189

$ws = $wb.Worksheets.Item('South Chicago')
$ws.Range("K2").Value = 33
$ws.Range("K7").Value = 58

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Range("K3").Value = 49
$ws.Range("K6").Value = 41
$ws.Range("K7").Value = 133

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Range("K3").Value = 13
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 54

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Range("K4").Value = 7
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 99

$ws = $wb.Worksheets.Item('New City')
$ws.Range("K2").Value = 19
$ws.Range("K7").Value = 88

$ws = $wb.Worksheets.Item('Fuller Park')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Range("K6").Value = 26
$ws.Range("K7").Value = 95
$ws.Range("K8").Value = 189
$ws.Range("K9").Value = 17
$ws.Range("K16").Value = 7
$ws.Range("K18").Value = 25
$ws.Range("K19").Value = 86
$ws.Range("K20").Value = 75
$ws.Range("K24").Value = 9
$ws.Range("K27").Value = 39
$ws.Range("K29").Value = 157
$ws.Range("K30").Value = 7
$ws.Range("K32").Value = 5
$ws.Range("K33").Value = 133
$ws.Range("J36").Value = 402
$ws.Range("K37").Value = 99
$ws.Range("K39").Value = 6
$ws.Range("K41").Value = 31
$ws.Range("K42").Value = 101
$ws.Range("K44").Value = 27
$ws.Range("K46").Value = 7
$ws.Range("K47").Value = 22
$ws.Range("K48").Value = 28
$ws.Range("K50").Value = 16
$ws.Range("K51").Value = 45
$ws.Range("K52").Value = 79
$ws.Range("K54").Value = 57
$ws.Range("K62").Value = 2
$ws.Range("K64").Value = 18
$ws.Range("K65").Value = 88
$ws.Range("K66").Value = 14
$ws.Range("K67").Value = 137
$ws.Range("K68").Value = 9
$ws.Range("K70").Value = 9
$ws.Range("K73").Value = 35
$ws.Range("K75").Value = 10
$ws.Range("K78").Value = 44
$ws.Range("K79").Value = 89
$ws.Range("K80").Value = 12
$ws.Range("K83").Value = 58
$ws.Range("K85").Value = 159
$ws.Range("K88").Value = 42
$ws.Range("K91").Value = 33
$ws.Range("K95").Value = 54
$ws.Range("K96").Value = 50
$ws.Range("J101").Value = 29250
$ws.Range("K101").Value = 3195

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Range("K2").Value = 39
$ws.Range("K7").Value = 137

$ws = $wb.Worksheets.Item('Loop')
$ws.Range("K3").Value = 19
$ws.Range("K6").Value = 22
$ws.Range("K7").Value = 57

$ws = $wb.Worksheets.Item('Englewood')
$ws.Range("K2").Value = 39
$ws.Range("K5").Value = 5
$ws.Range("K6").Value = 64
$ws.Range("K7").Value = 157

$ws = $wb.Worksheets.Item('Lake View')
$ws.Range("K2").Value = 8
$ws.Range("K6").Value = 10
$ws.Range("K7").Value = 28

$ws = $wb.Worksheets.Item('Chatham')
$ws.Range("K3").Value = 29
$ws.Range("K5").Value = 4
$ws.Range("K7").Value = 86

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Range("K6").Value = 13
$ws.Range("K7").Value = 27

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Range("K6").Value = 8
$ws.Range("K7").Value = 26

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Range("K2").Value = 9
$ws.Range("K7").Value = 31

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Range("K3").Value = 36
$ws.Range("K7").Value = 101

$ws = $wb.Worksheets.Item('Rogers Park')
$ws.Range("K2").Value = 17
$ws.Range("K3").Value = 11
$ws.Range("K7").Value = 44

$ws = $wb.Worksheets.Item('Dunning')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Range("K2").Value = 3
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Range("K2").Value = 18
$ws.Range("K6").Value = 25
$ws.Range("K7").Value = 50

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Range("K2").Value = 14
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 33

$ws = $wb.Worksheets.Item('Roseland')
$ws.Range("K2").Value = 30
$ws.Range("K3").Value = 30
$ws.Range("K7").Value = 89

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Range("K2").Value = 3
$ws.Range("K4").Value = 2
$ws.Range("K7").Value = 18

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Range("K2").Value = 16
$ws.Range("K3").Value = 22
$ws.Range("K6").Value = 34
$ws.Range("K7").Value = 75

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Range("K6").Value = 7
$ws.Range("K7").Value = 25

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Range("J6").Value = 120
$ws.Range("J7").Value = 402

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Range("K2").Value = 39
$ws.Range("K3").Value = 24
$ws.Range("K7").Value = 95

$ws = $wb.Worksheets.Item('Kenwood')
$ws.Range("K2").Value = 7
$ws.Range("K7").Value = 22

$ws = $wb.Worksheets.Item('Lincoln Square')
$ws.Range("K6").Value = 11
$ws.Range("K7").Value = 16

$ws = $wb.Worksheets.Item('Greektown')
$ws.Range("K5").Value = 4
$ws.Range("K6").Value = 6

$ws = $wb.Worksheets.Item('North Center')
$ws.Range("K3").Value = 3
$ws.Range("K7").Value = 14

$ws = $wb.Worksheets.Item('Avalon Park')
$ws.Range("K3").Value = 5
$ws.Range("K7").Value = 17

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Range("K2").Value = 11
$ws.Range("K7").Value = 35

$ws = $wb.Worksheets.Item('O''Hare')
$ws.Range("K2").Value = 2
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('United Center')
$ws.Range("K3").Value = 8
$ws.Range("K6").Value = 24
$ws.Range("K7").Value = 42

$ws = $wb.Worksheets.Item('Galewood')
$ws.Range("K7").Value = 5
$ws.Range("K4").Value = 1

$ws = $wb.Worksheets.Item('Edgewater')
$ws.Range("K6").Value = 18
$ws.Range("K7").Value = 39

$ws = $wb.Worksheets.Item('Pullman')
$ws.Range("K2").Value = 4
$ws.Range("K7").Value = 10

$ws = $wb.Worksheets.Item('Little Italy, UIC')
$ws.Range("K6").Value = 16
$ws.Range("K7").Value = 45

$ws = $wb.Worksheets.Item('North Park')
$ws.Range("K3").Value = 2
$ws.Range("K7").Value = 9

$ws = $wb.Worksheets.Item('South Shore')
$ws.Range("K2").Value = 61
$ws.Range("K4").Value = 10
$ws.Range("K6").Value = 37
$ws.Range("K7").Value = 159

$ws = $wb.Worksheets.Item('Rush & Division')
$ws.Range("K3").Value = 4
$ws.Range("K7").Value = 12

$ws = $wb.Worksheets.Item('Little Village')
$ws.Range("K6").Value = 39
$ws.Range("K7").Value = 79

$ws = $wb.Worksheets.Item('Bucktown')
$ws.Range("K6").Value = 6
$ws.Range("K7").Value = 7

$ws = $wb.Worksheets.Item('Museum Campus')
$ws.Range("K7").Value = 2
$ws.Range("K2").Value = 1
